$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted above the old row 116, pushing the
# existing records (old rows 116-134) down by one row (new rows 117-135).
$ws.Rows.Item(116).Insert()

# Populate the newly inserted row 116 with the new "Puerro" price record.
$ws.Range("A116").Value = 10
$ws.Range("B116").Value = "Vega Modelo de Temuco"
$ws.Range("C116").Value = "La Araucanía"
$ws.Range("D116").Value2 = 44474
$ws.Range("E116").Value = 9
$ws.Range("F116").Value = 100112005
$ws.Range("G116").Value = "Puerro"
$ws.Range("H116").Value = "Azul de Maquehue"
$ws.Range("I116").Value = "Primera"
$ws.Range("J116").Value = 30
$ws.Range("K116").Value = 7000
$ws.Range("L116").Value = 7000
$ws.Range("M116").Value = 7000
$ws.Range("N116").Value = "$/docena de paquetes"
$ws.Range("O116").Value = "Provincia de Cautín"
$ws.Range("P116").Value = 583
$ws.Range("Q116").Value = 12
$ws.Range("R116").Value = "Hortaliza"
